$wb = $excel.ActiveWorkbook

# --- Update the "Date" metadata property (sheet "Metadata", row 8: Date / value) ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2022-08-02T11:06:34-05:00"

# --- Remove the erroneous "CC0010583 / cyclophosphamide" row from the RxNorm include list ---
# (sheet "Include from RxNorm", row 5 currently holds A5=CC0010583, B5=cyclophosphamide,
# a duplicate of row 2's concept). Deleting the entire row shifts all following rows up by one.
$rxnorm = $wb.Worksheets.Item("Include from RxNorm")
$rxnorm.Range("A5:B5").EntireRow.Delete()
